# Auto-generated edit script: updates cryptos.xlsx price/volume data
# per commit 'Updated cryptos list on Thu Aug 15 18:32:02 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.904.90'
$ws.Range("E2").Value = '  -2.30%  '
$ws.Range("D3").Value = '2.556.73'
$ws.Range("E3").Value = '  -3.71%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'520.94"
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").Value = "'139.79"
$ws.Range("E6").Value = '  -3.58%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = "'0.560"
$ws.Range("E8").Value = '  -2.08%  '
$ws.Range("D9").Value = "'6.56"
$ws.Range("E9").Value = '  -6.87%  '
$ws.Range("D10").Value = "'0.0990"
$ws.Range("E10").Value = '  -3.83%  '
$ws.Range("D11").Value = "'0.324"
$ws.Range("E11").Value = '  -3.08%  '
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '3.006.40'
$ws.Range("E13").Value = '  -3.77%  '
$ws.Range("D14").Value = '57.930.23'
$ws.Range("E14").Value = '  -2.25%  '
$ws.Range("D15").Value = "'19.97"
$ws.Range("E15").Value = '  -5.52%  '
$ws.Range("D16").Value = '2.591.33'
$ws.Range("E16").Value = '  -2.12%  '
$ws.Range("E17").Value = '  -3.18%  '
$ws.Range("D18").Value = "'333.93"
$ws.Range("E18").Value = '  -2.08%  '
$ws.Range("D19").Value = "'4.28"
$ws.Range("E19").Value = '  -2.51%  '
$ws.Range("D20").Value = "'10.13"
$ws.Range("E20").Value = '  -2.52%  '
$ws.Range("D21").Value = "'6.12"
$ws.Range("E21").Value = '  -4.03%  '
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = "'64.96"
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("B24").Value = 'Kaspa'
$ws.Range("C24").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D24").Value = "'0.163"
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Value = "'0.401"
$ws.Range("E26").Value = '  -4.55%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.703.64'
$ws.Range("E27").Value = '  -2.55%  '
$ws.Range("D28").Value = "'6.94"
$ws.Range("E28").Value = '  -2.83%  '
$ws.Range("D29").Value = '0.0₃0752'
$ws.Range("E29").Value = '  -6.67%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = "'6.16"
$ws.Range("E31").Value = '  -7.86%  '
$ws.Range("E32").Value = '  -1.52%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").Value = "'18.42"
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("D35").Value = "'3.97"
$ws.Range("E35").Value = '  -4.70%  '
$ws.Range("D36").Value = "'1.12"
$ws.Range("E36").Value = '  -6.45%  '
$ws.Range("D37").Value = "'0.832"
$ws.Range("E37").Value = '  -7.14%  '
$ws.Range("D38").Value = "'35.90"
$ws.Range("E38").Value = '  -2.41%  '
$ws.Range("D39").Value = "'0.821"
$ws.Range("E39").Value = '  -6.89%  '
$ws.Range("D40").Value = "'1.40"
$ws.Range("E40").Value = '  -5.12%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = "'3.45"
$ws.Range("E42").Value = '  -4.41%  '
$ws.Range("D43").Value = "'0.0954"
$ws.Range("E43").Value = '  -1.89%  '
$ws.Range("D44").Value = "'10.62"
$ws.Range("E44").Value = '  +0.78%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = "'262.12"
$ws.Range("E45").Value = '  -4.80%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = "'0.577"
$ws.Range("E46").Value = '  -6.58%  '
$ws.Range("E47").Value = '  -3.19%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.985.73'
$ws.Range("E48").Value = '  -2.47%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'18.46"
$ws.Range("E49").Value = '  -8.27%  '
$ws.Range("E50").Value = '  -3.69%  '
$ws.Range("D51").Value = "'4.53"
$ws.Range("E51").Value = '  -5.86%  '
